# This revision of Presentation1.pptx is a build artifact from the
# PowerPoll Visual Studio project (bin/Debug/Presentation1.pptx). Diffing
# the canonical OOXML against the previous commit shows that:
#
#   1. Every <p:sldMasterId>/<p:sldId>/<p:sldLayoutId>/webextensionref/
#      blip r:id (and the matching r:embed/r:id values in the .rels
#      parts) were re-minted with new random hex suffixes. This is
#      simply what happens whenever the project is rebuilt and the
#      Open XML writer resaves the package - the relationship-id
#      strings are opaque, randomly generated identifiers with no
#      semantic meaning, and they are not something a user (or a
#      PowerPoint automation/VBA script) ever sets explicitly.
#   2. The only change that carries any real meaning is the identity
#      GUID of the embedded Office Add-in (the PowerPoll web extension)
#      stored in ppt/slides/udata/data.xml
#      (we:webextension/@id, "{23a6934e-...}" -> "{c8afd194-...}").
#      That GUID is minted by Visual Studio/PowerPoint when the add-in
#      manifest/snapshot is (re)inserted into the deck; it is part of
#      the webextension part (content type
#      application/vnd.ms-office.webextension+xml), which is a
#      completely different mechanism from CustomXMLParts and has
#      never been exposed on the PowerPoint Shape/Presentation object
#      model - there is no Shapes.AddWebExtension, no
#      Shape.WebExtension, and no supported way for a VBA/COM client
#      (or this host's equivalent automation surface) to read or set
#      it. It can only be produced by PowerPoint's own "Insert my
#      Add-in" UI flow / the add-in's own manifest tooling, never by
#      driving the object model.
#
# In other words, nothing in this commit is actually reachable through
# PowerPoint COM automation: the visible content of the single slide
# (title/subtitle placeholders, the OfficeApp picture fallback, the
# embedded image) is identical before and after, and the only textual
# change is to opaque identifiers that automation cannot touch. The
# faithful COM-interop reproduction is therefore to simply open the
# deck and leave its content exactly as authored, without introducing
# any spurious shape/text changes of our own.

$p = $ppt.ActivePresentation

# Touch the object model (read-only) so this script demonstrably runs
# against the live presentation without mutating any user-visible
# content - every shape/text/picture already matches the target state.
$s = $p.Slides.Item(1)
$null = $s.Shapes.Count
